# Update "Förändrad" (column C) dates from 2024-01-16 (45307) to 2024-01-17 (45308)
# for all existing data rows (2-27), then remove the last data row (row 28,
# "A 1668-2024") which no longer belongs in the sheet, and tidy up row 27's
# height so it reverts to the sheet's standard (non-custom) row height.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value = 45308
}

# Remove the now-obsolete last row (row 28, "A 1668-2024").
$ws.Rows("28:28").Delete()

# Row 27 previously had an explicit custom height; after the update it should
# use the sheet's standard height again (no ht/customHeight attributes).
$ws.Rows("27:27").AutoFit()
